$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B values for rows 2-3 (introduce new shared strings rx0, tx1)
$ws.Range("B2").Value = "rx0"
$ws.Range("B3").Value = "tx1"

# Column B numeric values for rows 4-13 (plain numbers, no new shared strings)
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 8
$ws.Range("B11").Value = 9
$ws.Range("B12").Value = 10
$ws.Range("B13").Value = 11

# Column B numeric values for rows 14-15 (plain numbers, no new shared strings)
$ws.Range("B14").Value = 12
$ws.Range("B15").Value = 13

# Column B for rows 16-19 (introduce new shared strings tx3, rx3, tx2, rx2)
$ws.Range("B16").Value = "tx3"
$ws.Range("B17").Value = "rx3"
$ws.Range("B18").Value = "tx2"
$ws.Range("B19").Value = "rx2"

# Column B for row 20 (reuses tx1)
$ws.Range("B20").Value = "tx1"

# New row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "rx1"

# New row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "sda"

# New row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "scl"

# Column C for rows 23, 22, 18, 19 (introduce SCL(LCD), SDA(LCD), rasp pi (rx), rasp pi (tx))
$ws.Range("C23").Value = "SCL(LCD)"
$ws.Range("C22").Value = "SDA(LCD)"
$ws.Range("C18").Value = "rasp pi (rx)"
$ws.Range("C19").Value = "rasp pi (tx)"

# Column C for rows 14, 15 (introduce trig back sonar, echo back sonar)
$ws.Range("C14").Value = "trig back sonar "
$ws.Range("C15").Value = "echo back sonar"

# View changes: topLeftCell A7, selection D15
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D15").Select()
